# Add a new "dbsource" one-hot-encoding block to the factorization table.
# This mirrors the author's classification.py change that added
# one-hot-encoding for the "dbsource" feature (carevue / both / metavision).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right above the current "ethnicity" block (row 5),
# pushing everything else down and keeping the existing number formats
# (row styles) intact.
$ws.Rows("5:7").Insert()

# Populate the newly inserted rows with the dbsource factorization data.
$ws.Range("A5").Value = "dbsource"
$ws.Range("B5").Value = "carevue"
$ws.Range("C5").Value = -1

$ws.Range("A6").Value = "dbsource"
$ws.Range("B6").Value = "both"
$ws.Range("C6").Value = 0

$ws.Range("A7").Value = "dbsource"
$ws.Range("B7").Value = "metavision"
$ws.Range("C7").Value = 1

# Match the author's final selected cell.
[void]$ws.Range("C8").Select()
